$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")
$ws.Activate() | Out-Null

$ws.Range("B6").Value = 55.5
$ws.Range("C6").Value = 55.5
$ws.Range("D6").Value = 58
$ws.Range("E6").Value = 51.5
$ws.Range("F6").Value = 53.5
$ws.Range("G6").Value = 55.5

$ws.Range("B7").Value = 53.5
$ws.Range("C7").Value = 53.5
$ws.Range("D7").Value = 55
$ws.Range("E7").Value = 47
$ws.Range("F7").Value = 53.5
$ws.Range("G7").Value = 57.5

$ws.Range("B8").Value = 54
$ws.Range("C8").Value = 55.5
$ws.Range("D8").Value = 56.5
$ws.Range("E8").Value = 48.5
$ws.Range("F8").Value = 55
$ws.Range("G8").Value = 48.5

$ws.Range("B8").Select() | Out-Null
